$d = $word.ActiveDocument

# The document currently holds a single, empty paragraph whose paragraph
# mark already carries the intended run formatting (Lato font, bold,
# color 2D3B45, black underline color). The edit adds one run of text,
# "fsa", using that same formatting.

$sel = $word.Selection
$sel.EndKey(6) | Out-Null       # wdStory -> move to the very end of the document
$sel.TypeText("fsa") | Out-Null # insert the new run text

# Re-apply the formatting explicitly so the inserted run carries it
# (ascii/eastAsia/hAnsi/cs all set to Lato, bold incl. complex-script,
# font color 2D3B45, underline color black).
$sel.Font.Name = "Lato"
$sel.Font.NameAscii = "Lato"
$sel.Font.NameFarEast = "Lato"
$sel.Font.NameOther = "Lato"
$sel.Font.NameBi = "Lato"
$sel.Font.Bold = $true
$sel.Font.BoldBi = $true
$sel.Font.Color = 4537133     # 0x45 3B 2D (BGR) == RGB 2D3B45
$sel.Font.UnderlineColor = 0  # black (000000)
